# Insert a new data row at row 82 (shifting the existing rows 82-136 down
# to 83-137) and populate the new row with the latest Mango price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(82).Insert()

$ws.Range("A82").Value = 7
$ws.Range("B82").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C82").Value = "Ñuble"
$ws.Range("D82").Value = 45068
$ws.Range("E82").Value = 16
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100108
$ws.Range("H82").Value = "Tropicales y subtropicales"
$ws.Range("I82").Value = 100108002
$ws.Range("J82").Value = "Mango"
$ws.Range("K82").Value = "Sin especificar"
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 60
$ws.Range("N82").Value = 8000
$ws.Range("O82").Value = 8000
$ws.Range("P82").Value = 8000
$ws.Range("Q82").Value = "$/bandeja 4 kilos"
$ws.Range("R82").Value = "Perú"
$ws.Range("S82").Value = 2000
$ws.Range("T82").Value = 4
